# Linter fix: replace underscore with a space in the "birth_date" column
# header on every sheet that contains it (Animal, NamedThing, Person,
# Animal1, NamedThing1, Person1).

$wb = $excel.ActiveWorkbook

$targets = @{
    "Animal"        = "C1"
    "NamedThing"     = "C1"
    "Person"         = "B1"
    "Animal1"        = "C1"
    "NamedThing1"    = "C1"
    "Person1"        = "B1"
}

foreach ($sheetName in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellRef = $targets[$sheetName]
    $cell = $ws.Range($cellRef)
    if ($cell.Value2 -eq "birth_date") {
        $cell.Value = "birth date"
    }
}
